# Matriz Casos Prueba V1.0
# Renumbers the "Caso de Prueba N" headings and their matching "Identificador"
# (CP-N) values into sequential document order, renames the sheet, and
# updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---
$ws.Name = "Hoja 1"

# --- Renumber "Caso de Prueba N: ..." titles (column A) and their
#     corresponding "Identificador" -> CP-N values (column B) so that the
#     whole matrix is numbered sequentially 1-8 from top to bottom. ---

# Inicio de Sesion block (CP-7 / CP-8 -> CP-1 / CP-2, titles unchanged)
$ws.Range("B7").Value  = "CP-1"
$ws.Range("B18").Value = "CP-2"

# Cambio de Contrasena block (was "Caso de Prueba 5/6", CP-5/CP-6 -> 3/4)
$ws.Range("A29").Value = "Caso de Prueba 3: Actualizar contraseña correctamente"
$ws.Range("B30").Value = "CP-3"

$ws.Range("A44").Value = "Caso de Prueba 4: Actualizar contraseña incorrectamente"
$ws.Range("B45").Value = "CP-4"

# Gestion de Sugerencias block (was "Caso de Prueba 1/2", CP-1/CP-2 -> 5/6)
$ws.Range("A60").Value = "Caso de Prueba 5: Enviar sugerencias exitosamente"
$ws.Range("B61").Value = "CP-5"

$ws.Range("A76").Value = "Caso de Prueba 6: Intentar enviar sugerencias sin activar el check"
$ws.Range("B77").Value = "CP-6"

# Actualizacion de Datos Personales block (was "Caso de Prueba 3/4", CP-3/CP-4 -> 7/8)
$ws.Range("A93").Value = "Caso de Prueba 7: Actualizar datos personales correctamente"
$ws.Range("B94").Value = "CP-7"

$ws.Range("A108").Value = "Caso de Prueba 8: Actualizar datos personales incorrectamente"
$ws.Range("B109").Value = "CP-8"

# --- Update the active window scroll position / selection ---
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 144
    $win.ScrollColumn = 1
} catch {
}
$ws.Range("A108:B108").Select()
